$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1052.745727539062,
    12407.5595703125,
    21590.033203125,
    40462.88671875,
    47262.29296875,
    55567.96484375,
    65615.4296875,
    89607.7265625,
    101289.3984375,
    100351.875,
    71273.6640625,
    39045.24609375,
    6526.80712890625,
    63760.02734375,
    2031.147094726562,
    37300.88671875,
    69052.6484375,
    55886.10546875,
    45507.74609375,
    73801.9296875,
    5810.43212890625,
    62747.57421875,
    81431.2890625,
    84421.6484375,
    63407.74609375,
    3075.182373046875,
    61625.43359375,
    44256.55859375,
    29287.939453125,
    59647.46484375,
    55863.35546875,
    86938.859375,
    57020.62109375,
    109347.4921875,
    81926.2265625,
    2788.225341796875,
    39394.50390625,
    55592.93359375,
    68053.2265625,
    87148.3203125,
    133978.6875,
    2564.358154296875,
    58986.44921875,
    20910.427734375,
    51889.21484375,
    46338.16796875,
    66783.4140625,
    71530.6171875,
    55231.33984375,
    68749.8359375,
    87173.4921875,
    108401.890625,
    113940.140625,
    124625.3203125,
    87998.1015625,
    50575.08984375,
    51933.91796875,
    53069.02734375,
    80529.7734375,
    75914.6484375,
    92329.3515625,
    102383.84375,
    113539.8125,
    126727.9609375,
    62825.69921875,
    63304.60546875,
    78363.0078125,
    96411.46875,
    55690.60546875
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $values[$i]
}

Write-Host "Updated $($values.Length) cells, rows $startRow to $($startRow + $values.Length - 1)"
